# Commit: "#5: property boat&car done"
#
# The 汽車 (car) sheet (sheet3) is rewritten to follow the common
# property-table template used by the other sheets (name / capacity /
# owner / register_date / register_reason / acquire_value /
# property_category / category / date / legislator_name / legislator_id /
# source_file / index), and the data row is filled in accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# ---- Header row (row 1) -------------------------------------------------
# Existing header cells B1:G1 get new text; H1:N1 are brand-new cells.
$ws.Range("B1").Value2 = "name"
$ws.Range("C1").Value2 = "capacity"
$ws.Range("D1").Value2 = "owner"
$ws.Range("E1").Value2 = "register_date"
$ws.Range("F1").Value2 = "register_reason"
$ws.Range("G1").Value2 = "acquire_value"
$ws.Range("H1").Value2 = "property_category"
$ws.Range("I1").Value2 = "category"
$ws.Range("J1").Value2 = "date"
$ws.Range("K1").Value2 = "legislator_name"
$ws.Range("L1").Value2 = "legislator_id"
$ws.Range("M1").Value2 = "source_file"
$ws.Range("N1").Value2 = "index"

# New header cells need the same formatting (bold/border) as the rest of
# row 1, so copy the format from an existing header cell onto them.
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---- Data row (row 2) ----------------------------------------------------
# B2 keeps the same text ("TOYOTA(自用小客車）") - just rewritten so the
# sheet is self-consistent after the shared-string shuffle.
$ws.Range("B2").Value2 = "TOYOTA(自用小客車）"
$ws.Range("H2").Value2 = "land"
$ws.Range("I2").Value2 = "normal"
$ws.Range("J2").Value2 = "2012-02-01"
$ws.Range("K2").Value2 = "馬文君"
$ws.Range("L2").Value2 = 1724
$ws.Range("M2").Value2 = "tmpb9501"
$ws.Range("N2").Value2 = 42

# New data cells need the same formatting as the rest of row 2.
$ws.Range("B2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
